# feat: implement basic whatsapp bot ping pong
# Mark the three "WhatsApp Bot" tracker rows (Setup project Node.js,
# Install whatsapp-web.js, Implementasi bot dasar (ping-pong)) as done,
# and leave the active selection on E13.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("F10").Value = $true
$ws.Range("F11").Value = $true
$ws.Range("F12").Value = $true

$ws.Range("E13").Select()
